$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data (row 11) -- order matches shared-string append order
$ws.Range("E11").Value = "https://leetcode.com/problems/average-time-of-process-per-machine/solutions/3537533/detailed-explanation-of-join-round-avg-group-by/?envType=study-plan-v2&envId=top-sql-50 "
$ws.Range("D11").Value = "We can solve with JOIN, ROUND, AVG, GROUP BY. First we need to join the same table to have starting time and ending time in the same line. Then calculate the difference between timestamps and group the values by machine id."
$ws.Range("A11").Value = "1661. Average Time of Process per Machine"
$ws.Range("B11").Value = "Easy"
$ws.Range("C11").Value = "Basic Joins"

# Match the formatting used by the other data rows
$ws.Range("B11").Interior.Color = 5287936

# Add hyperlink for the Link column, matching the other rows
$ws.Hyperlinks.Add($ws.Range("E11"), "https://leetcode.com/problems/average-time-of-process-per-machine/solutions/3537533/detailed-explanation-of-join-round-avg-group-by/?envType=study-plan-v2&envId=top-sql-50 ") | Out-Null

# Hyperlinks.Add mints its own cell-style xf; re-apply the exact format used
# by the other Link-column cells so E11 matches them byte-for-byte.
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Grow the table to include the new row
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("A1:E11"))

# Update the active selection like in the authored workbook
$ws.Range("E17").Select()
